$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D-column price cells whose new values would
# otherwise be auto-parsed as numbers by Excel (a value with at most one
# "." looks numeric; a value with two "." like "66.139.46" stays text on
# its own). Setting NumberFormat to Text ("@") before assigning keeps
# these as literal strings, matching the source data which is textual.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.139.46"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "3.563.32"
$ws.Range("E3").Value = "  +2.19%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "605.84"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "144.44"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "3.562.63"
$ws.Range("E7").Value = "  +2.21%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "4.169.09"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Value = "30.07"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "3.572.63"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "66.213.21"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "11.45"
$ws.Range("E19").Value = "  +6.36%  "
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("D22").Value = "431.07"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("D24").Value = "79.78"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("D25").Value = "3.707.73"
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("D32").Value = "3.562.21"
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  -7.88%  "
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("D39").Value = "5.55"
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").Value = "173.86"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").Value = "5.18"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").Value = "0.893"
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").Value = "1.94"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").Value = "45.97"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").Value = "25.01"
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "22.96"
$ws.Range("E51").Value = "  +4.99%  "
